$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fill in new run data (Lauf/Schritte = column B, Distanz = column C) for each problem block ---
# Block 1: rows 4-23
$bVals = @(3,2,2,1,2,2,2,5,1,1,1,0,1,1,3,1,5,1,3,1)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(4 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(4 + $i, 3).Value = $cVals[$i]
}

# Block 2: rows 29-48
$bVals = @(25,21,13,21,16,15,15,13,21,17,17,29,17,24,24,20,15,26,19,7)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(29 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(29 + $i, 3).Value = $cVals[$i]
}

# Block 3: rows 54-73
$bVals = @(18,18,27,28,42,16,29,22,38,28,32,27,24,24,19,34,26,19,16,17)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(54 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(54 + $i, 3).Value = $cVals[$i]
}

# Block 4: rows 79-98
$bVals = @(20,27,6,5,300,10,23,5,26,0,21,14,300,11,6,300,11,12,15,16)
$cVals = @(0,0,0,0,1,0,0,5,0,1,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(79 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(79 + $i, 3).Value = $cVals[$i]
}

# Block 5: rows 104-123
$bVals = @(1,2,4,2,1,3,3,1,3,2,9,1,5,3,4,3,1,3,2,0)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(104 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(104 + $i, 3).Value = $cVals[$i]
}

# Block 6: rows 129-148
$bVals = @(4,8,8,3,4,2,4,6,4,2,2,4,5,12,4,1,6,4,4,1)
$cVals = @(1,2,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(129 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(129 + $i, 3).Value = $cVals[$i]
}

# Block 7: rows 154-173
$bVals = @(11,20,22,16,18,7,9,12,7,11,12,14,12,19,10,9,11,7,9,9)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,2,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(154 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(154 + $i, 3).Value = $cVals[$i]
}

# Block 8: rows 179-198
$bVals = @(8,11,8,7,14,12,5,8,5,13,4,5,12,16,14,15,7,14,9,16)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(179 + $i, 2).Value = $bVals[$i]
    $ws.Cells.Item(179 + $i, 3).Value = $cVals[$i]
}

# --- Average rows: keep column B AVERAGE formula, remove the column C AVERAGE (no longer used) ---
$ws.Range('B24').Formula = '=AVERAGE(B4:B23)'
$ws.Range('C24').ClearContents()
$ws.Range('B49').Formula = '=AVERAGE(B29:B48)'
$ws.Range('C49').ClearContents()
$ws.Range('B74').Formula = '=AVERAGE(B54:B73)'
$ws.Range('C74').ClearContents()
$ws.Range('B99').Formula = '=AVERAGE(B79:B98)'
$ws.Range('C99').ClearContents()
$ws.Range('B124').Formula = '=AVERAGE(B104:B123)'
$ws.Range('C124').ClearContents()
$ws.Range('B149').Formula = '=AVERAGE(B129:B148)'
$ws.Range('C149').ClearContents()
$ws.Range('B174').Formula = '=AVERAGE(B154:B173)'
$ws.Range('C174').ClearContents()
$ws.Range('B199').Formula = '=AVERAGE(B179:B198)'
$ws.Range('C199').ClearContents()

# --- Histogram COUNTIF formulas now reference the extended range C4:C224 ---
$ws.Range('G4').Formula = '=COUNTIF($C$4:$C$224,0)'
$ws.Range('G5').Formula = '=COUNTIF($C$4:$C$224,1)'
$ws.Range('G6').Formula = '=COUNTIF($C$4:$C$224,2)'
$ws.Range('G7').Formula = '=COUNTIF($C$4:$C$224,3)'
$ws.Range('G8').Formula = '=COUNTIF($C$4:$C$224,4)'
$ws.Range('G9').Formula = '=COUNTIF($C$4:$C$224,5)'
$ws.Range('G10').Formula = '=COUNTIF($C$4:$C$224,6)'
$ws.Range('G11').Formula = '=COUNTIF($C$4:$C$224,7)'
$ws.Range('G12').Formula = '=COUNTIF($C$4:$C$224,8)'
$ws.Range('G13').Formula = '=COUNTIF($C$4:$C$224,9)'
$ws.Range('G14').Formula = '=COUNTIF($C$4:$C$224,10)'

# --- Probability formula recovers from #DIV/0! once G4:G14 are populated ---
$ws.Range('F20').Formula = '=SUM(G4,G5)/SUM(G4:G14)'

# --- Restore the selection state recorded in the saved workbook ---
[void]$ws.Range('I23').Select()

$wb.Application.Calculate()
